$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    @(2, 10.89835490497452),
    @(3, 24.62750397673582),
    @(4, 20.88289491432329),
    @(5, 22.1714314217609),
    @(6, 22.2459716575656),
    @(7, 19.8627667704732),
    @(8, 17.51806114281132),
    @(9, 19.88765161502965),
    @(10, 13.87313780430065),
    @(11, 20.1905144245213),
    @(12, 22.20159577114279),
    @(13, 14.61544706729242),
    @(14, 14.2277478471413),
    @(15, 16.45236482271528),
    @(16, 14.37093776810019),
    @(17, 13.49500675559565),
    @(18, 18.63777464806455),
    @(19, 10.57455121237427),
    @(20, 14.51716305219858),
    @(21, 12.50492244635581),
    @(22, 11.99880493412905),
    @(23, 13.74415421457692),
    @(24, 11.30057257976017),
    @(25, 11.73164096260487),
    @(26, 8.036845902372278),
    @(27, 9.994628589008698),
    @(28, 13.53282529787364),
    @(29, 8.377997630179237),
    @(30, 8.432203183078542),
    @(31, 5.461042980876613),
    @(32, 8.86906268077874),
    @(33, 9.915064404167424),
    @(34, 10.39602741959368),
    @(35, 11.28706645135679),
    @(36, 8.15975534162385),
    @(37, 7.973164742951411),
    @(38, 7.086360762773353),
    @(39, 7.841572175589988),
    @(40, 5.981930266154478),
    @(41, 5.791030265471221),
    @(42, 6.26027017159663),
    @(43, 9.456097019158648),
    @(44, 8.109963575079547),
    @(45, 10.88650529577026),
    @(46, 12.40542978044638),
    @(47, 8.606416124972782),
    @(48, 8.9217683513466),
    @(49, 6.944412901889081),
    @(50, 8.950332879895939),
    @(51, 6.932683172820759),
    @(52, 8.346943068934337),
    @(53, 9.373061270529092),
    @(54, 6.675810945434108),
    @(55, 7.201672055830755),
    @(56, 7.462466999753644),
    @(57, 9.534936574661771),
    @(58, 7.957001635313219),
    @(59, 7.69888144684748),
    @(60, 8.124326528385666),
    @(61, 7.219530104558743),
    @(62, 6.566454011678218),
    @(63, 4.391146095249752),
    @(64, 3.93816971401219),
    @(65, 8.456032965627941),
    @(66, 4.234850931000324),
    @(67, 7.087035791471692),
    @(68, 4.752798452809401),
    @(69, 4.227564117610541),
    @(70, 6.615625445471693),
    @(71, 8.602491958863084),
    @(72, 5.223301369830125),
    @(73, 7.229241319596071),
    @(74, 4.162270926113592),
    @(75, 9.234831465217525),
    @(76, 8.061591374610117),
    @(77, 7.853951870358628),
    @(78, 6.939010895201619),
    @(79, 7.670603352320427),
    @(80, 9.455493754549593),
    @(81, 7.257951807573249),
    @(82, 9.421277742469471),
    @(83, 7.69714408247188),
    @(84, 7.13286838164791),
    @(85, 9.765923334843393),
    @(86, 7.165612879813466)
)

foreach ($pair in $values) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($row, 1).Value = $val
}

Write-Output "Updated $($values.Count) cells in column A"